$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    43.63717779732323,25.26993961287565,25.7139452826859,21.61282162320063,21.30751531249774,20.90322275874323,22.48756082702388,18.29439548919828,24.23377528778366,16.11400757820931,
    10.90894950325495,18.12865209125238,18.67173236435809,14.68173723046135,15.93782430269817,15.97868083297291,9.557168634852786,17.57070945552513,13.70835898906233,15.11769696404329,
    14.4545245314518,12.00596857699631,14.62082409183182,13.07653650146946,17.04084803783749,15.36164955716816,9.706612241580984,15.25898608711805,14.37294801123835,17.11357208974596,
    14.19681811435214,13.28797829480777,12.21367515580084,10.43522857299075,13.54590686294839,13.36611842234905,14.74325740490616,13.58740068732033,14.97260439818496,15.1488091706025,
    16.11641364992687,10.36442283400939,13.36202965610286,9.367137243731662,7.275907842615055,11.97103489066868,10.34613999164296,12.49896050983048,10.77407049346436,12.14767845824653,
    10.62567925887222,8.991245458666782,11.9275206734506,11.46363991443434,11.45638472121456,8.714185167188219,10.46915209224312,9.80045994206526,9.297621732168665,10.1045564530618,
    11.99008487396756,15.08287799329577,14.84994778659845,9.408075640176207,14.4158133168269,10.94148508299716,13.38614594368946,14.85946567691244,11.80260435244972,9.558591977764905,
    12.74895477290306,10.60342154405438,14.13817397854444,8.365922088962321,9.663712059897449,9.835720709566374,11.02575980958699,10.8767219896921,9.953928039987233,12.69608608807522,
    9.758379732171761,11.59011565465735,12.14361763526915,9.116733211906677,11.41770121112199
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
